# Shorten the "family business" description on the Kokemus (FI) and
# en-kokemus (EN) sheets by removing the trailing sentence about the
# restaurant's website, leave the cursor on the edited cell, bump the
# Taidot sheet zoom, and finish with en-kokemus as the active sheet.

$wb = $excel.ActiveWorkbook

$wsKokemus = $wb.Worksheets.Item("Kokemus")
$wsKokemus.Range("C2").Value = "Perheyritys jossa olen työskennellyt ravintolan keittiössä ja auttanut johtoa IT- asioiden kanssa. "
$null = $wsKokemus.Range("C2").Select()

$wsEnKokemus = $wb.Worksheets.Item("en-kokemus")
$wsEnKokemus.Range("C2").Value = "A family run business where I have worked in the restaurant's kitchen and helped magement with IT- issues."
$null = $wsEnKokemus.Range("C2").Select()

$wsTaidot = $wb.Worksheets.Item("Taidot")
$wsTaidot.Activate()
$excel.ActiveWindow.Zoom = 225

$wsEnKokemus.Activate()
